$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 55, pushing existing rows 55-90 down to 56-91.
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new record.
$ws.Cells.Item(55, 1).Value = 1
$ws.Cells.Item(55, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(55, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(55, 4).Value = 44767
$ws.Cells.Item(55, 5).Value = 15
$ws.Cells.Item(55, 6).Value = 100112021
$ws.Cells.Item(55, 7).Value = "Ají"
$ws.Cells.Item(55, 8).Value = "Inferno"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 140
$ws.Cells.Item(55, 11).Value = 11000
$ws.Cells.Item(55, 12).Value = 12000
$ws.Cells.Item(55, 13).Value = 11500
$ws.Cells.Item(55, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(55, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(55, 16).Value = 767
$ws.Cells.Item(55, 17).Value = 15
$ws.Cells.Item(55, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item(55, 4).NumberFormat = $ws.Cells.Item(56, 4).NumberFormat
